# The workbook tracks a weekly price list. A new week's reading was
# inserted at the top of the data block (row 23), pushing the existing
# rows 23-69 down to 24-70. Reproduce that by inserting a new row at 23
# and then populating it: the four "identity" columns (market, region,
# category info, unit, origin, kg) are copied from the row that used to
# occupy row 23 (now at row 24, unchanged), while the week-specific
# figures (date, volume, min/max/weighted price, price per kg) get the
# new values from the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 23; rows 23:69 shift down to 24:70.
$ws.Rows.Item(23).Insert()

# Columns that stay constant for this market/category combo - copy them
# from the row directly below (which held this data before the insert).
$ws.Cells.Item(23, 1).Value = $ws.Cells.Item(24, 1).Value()
$ws.Cells.Item(23, 2).Value = $ws.Cells.Item(24, 2).Value()
$ws.Cells.Item(23, 3).Value = $ws.Cells.Item(24, 3).Value()
$ws.Cells.Item(23, 5).Value = $ws.Cells.Item(24, 5).Value()
$ws.Cells.Item(23, 6).Value = $ws.Cells.Item(24, 6).Value()
$ws.Cells.Item(23, 7).Value = $ws.Cells.Item(24, 7).Value()
$ws.Cells.Item(23, 8).Value = $ws.Cells.Item(24, 8).Value()
$ws.Cells.Item(23, 9).Value = $ws.Cells.Item(24, 9).Value()
$ws.Cells.Item(23, 14).Value = $ws.Cells.Item(24, 14).Value()
$ws.Cells.Item(23, 15).Value = $ws.Cells.Item(24, 15).Value()
$ws.Cells.Item(23, 17).Value = $ws.Cells.Item(24, 17).Value()
$ws.Cells.Item(23, 18).Value = $ws.Cells.Item(24, 18).Value()

# New weekly reading values.
$ws.Cells.Item(23, 4).Value = 44672    # D: Fecha
$ws.Cells.Item(23, 10).Value = 220     # J: Volumen
$ws.Cells.Item(23, 11).Value = 10000   # K: Precio minimo
$ws.Cells.Item(23, 12).Value = 11000   # L: Precio maximo
$ws.Cells.Item(23, 13).Value = 10455   # M: Precio promedio ponderado
$ws.Cells.Item(23, 16).Value = 1046    # P: Precio $/Kg
